# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with a newer snapshot, cell by cell, exactly as the data
# refresh workflow does.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.152.91"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.556.03"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'605.45"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'144.17"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "3.554.86"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'7.81"
$ws.Range("E11").Value = "  -3.23%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "4.156.20"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'30.01"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "3.549.76"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "66.191.50"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'11.32"
$ws.Range("E19").Value = "  +6.02%  "
$ws.Range("D20").Value = "'6.19"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").Value = "'429.99"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").Value = "'79.72"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("D25").Value = "3.698.26"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "3.551.67"
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "'25.39"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("E35").Value = "  -9.21%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "'174.41"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "'0.0844"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "'5.18"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'0.887"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'46.08"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'24.75"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'7.11"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'22.95"
$ws.Range("E51").Value = "  +3.06%  "
